# Update odds values for rows 2 and 5 (Jogos da Semana FlashScore 2025-02-20)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 2.5
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 1.91
$ws.Range("L2").Value = 3.4
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 2.38
$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 1.8
$ws.Range("S2").Value = 2.7
$ws.Range("T2").Value = 1.44
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 1.14
$ws.Range("Y2").Value = 1.62
$ws.Range("Z2").Value = 2.2
$ws.Range("AA2").Value = 2.2
$ws.Range("AB2").Value = 1.62
$ws.Range("AC2").Value = 7
$ws.Range("AD2").Value = 13
$ws.Range("AF2").Value = 34
$ws.Range("AI2").Value = 6
$ws.Range("AK2").Value = 19
$ws.Range("AO2").Value = 11
$ws.Range("AP2").Value = 23
$ws.Range("AQ2").Value = 26

# Row 5
$ws.Range("G5").Value = 2.77
$ws.Range("I5").Value = 2.37
$ws.Range("J5").Value = 3.35
$ws.Range("L5").Value = 2.9
$ws.Range("T5").Value = 1.88
$ws.Range("W5").Value = 2.67
$ws.Range("AB5").Value = 2.1
$ws.Range("AC5").Value = 9.75
$ws.Range("AD5").Value = 15
$ws.Range("AE5").Value = 10
$ws.Range("AF5").Value = 35
$ws.Range("AG5").Value = 23
$ws.Range("AH5").Value = 28
$ws.Range("AL5").Value = 50
$ws.Range("AM5").Value = 9.5
$ws.Range("AO5").Value = 9
$ws.Range("AQ5").Value = 18
$ws.Range("AR5").Value = 24

$wb.Save()
